$wb = $excel.ActiveWorkbook

# ---- Layer0 sheet ----
$ws1 = $wb.Worksheets.Item("Layer0")

$ws1.Range("B2").Value = 2.992399742986195
$ws1.Range("C2").Value = -3.115361607584886
$ws1.Range("D2").Value = 6.509189091274481
$ws1.Range("E2").Value = -0.9074571694696207

$ws1.Range("B3").Value = -0.02530588994318195
$ws1.Range("C3").Value = -0.9781551949651406
$ws1.Range("D3").Value = 2.463610012120194
$ws1.Range("E3").Value = 3.023334633946135

$ws1.Range("B4").Value = 3.693618807278793
$ws1.Range("C4").Value = 1.995690265495933
$ws1.Range("D4").Value = -8.604889775144191
$ws1.Range("E4").Value = 1.106936683455789

$ws1.Range("B5").Value = -2.016162561158827
$ws1.Range("C5").Value = 2.555806862681343
$ws1.Range("D5").Value = 0.9346310724968041
$ws1.Range("E5").Value = -2.133463750928092

$ws1.Range("B6").Value = 2.526680110979042
$ws1.Range("C6").Value = 2.229509010943596
$ws1.Range("D6").Value = 4.379522965084463
$ws1.Range("E6").Value = -3.468790923535672

$ws1.Range("B7").Value = 2.286771883344559
$ws1.Range("C7").Value = -8.545326068087178
$ws1.Range("D7").Value = 1.394992104309901
$ws1.Range("E7").Value = 0.05783737800959484

# ---- Layer1 sheet ----
$ws2 = $wb.Worksheets.Item("Layer1")

$ws2.Range("B2").Value = -4.390035857830386
$ws2.Range("C2").Value = -5.985822531774705
$ws2.Range("D2").Value = -3.737605885524056
$ws2.Range("E2").Value = -9.865883052535109
$ws2.Range("F2").Value = -4.598259568341395

$ws2.Range("B3").Value = 6.331951722016104
$ws2.Range("C3").Value = -9.54714353361236
$ws2.Range("D3").Value = -0.3147100914771699
$ws2.Range("E3").Value = -3.687887333112108
$ws2.Range("F3").Value = 7.473668636143406

$ws2.Range("B4").Value = -9.054268163880936
$ws2.Range("C4").Value = 5.224444448816967
$ws2.Range("D4").Value = -3.86142489380658
$ws2.Range("E4").Value = -5.718326167797028
$ws2.Range("F4").Value = 8.125615136262052

$ws2.Range("B5").Value = 3.630457069518193
$ws2.Range("C5").Value = 6.79068419089165
$ws2.Range("D5").Value = -9.646607528009406
$ws2.Range("E5").Value = 8.407713907826864
$ws2.Range("F5").Value = -8.216308862190841

$ws2.Range("B6").Value = -11.6578995476367
$ws2.Range("C6").Value = -6.476629873826448
$ws2.Range("D6").Value = 7.934443770723268
$ws2.Range("E6").Value = 9.748175293126748
$ws2.Range("F6").Value = -8.236350070575366
